# l26 600 qpsk wip
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LoopFilter LPF")
$ws2 = $wb.Worksheets.Item("Branch LPF")

# --- Sheet "Branch LPF" ---
# Update scale bits B18: 1 -> 4
$ws2.Range("B18").Value = 4

# Update selection to B18
$ws2.Activate()
$ws2.Range("B18").Select()

# --- Sheet "LoopFilter LPF" ---
# Update scale bits B16: 8 -> 2
$ws1.Range("B16").Value = 2

# Update selection to C18 and keep this sheet as the active/selected tab
$ws1.Activate()
$ws1.Range("C18").Select()

$wb.Save()
